$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.664768333333334
$ws.Range("H2").Value = 19.994305
$ws.Range("I2").Value = 0.06516174319532789
$ws.Range("J2").Value = 0.0651617431953279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3103893333333334
$ws.Range("N2").Value = 0.931168
$ws.Range("Q2").Value = 2.068672999804445
$ws.Range("R2").Value = 18.61805699824
$ws.Range("S2").Value = 0.06516174319532789
$ws.Range("T2").Value = 0.0651617431953279

# Row 3 (FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 47.25592399999999
$ws.Range("H3").Value = 141.767772
$ws.Range("I3").Value = 0.4620233187619072
$ws.Range("J3").Value = 0.4620233187619072
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3103893333333334
$ws.Range("N3").Value = 0.931168
$ws.Range("Q3").Value = 14.66773474641067
$ws.Range("R3").Value = 132.009612717696
$ws.Range("S3").Value = 0.4620233187619072
$ws.Range("T3").Value = 0.4620233187619072

# Row 4 (M2)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.98736333333333
$ws.Range("H4").Value = 62.96209
$ws.Range("I4").Value = 0.2051944060881897
$ws.Range("J4").Value = 0.2051944060881898
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3103893333333334
$ws.Range("N4").Value = 0.931168
$ws.Range("Q4").Value = 6.514253713457778
$ws.Range("R4").Value = 58.62828342112
$ws.Range("S4").Value = 0.2051944060881897
$ws.Range("T4").Value = 0.2051944060881898

# Row 5 (sCs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.37233166666667
$ws.Range("H5").Value = 82.116995
$ws.Range("I5").Value = 0.2676205319545753
$ws.Range("J5").Value = 0.2676205319545753
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3103893333333334
$ws.Range("N5").Value = 0.931168
$ws.Range("Q5").Value = 8.496079777795556
$ws.Range("R5").Value = 76.46471800016
$ws.Range("S5").Value = 0.2676205319545753
$ws.Range("T5").Value = 0.2676205319545753
